# Apply updated "想去人数" (want-to-go count) figures to the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Row -> new value mapping for sheet "展览"
$updates1 = @{
    2  = 14
    3  = 1343
    6  = 60
    8  = 11612
    9  = 4370
    11 = 34
    15 = 1091
    17 = 33
    18 = 5058
    21 = 11327
    22 = 11250
    23 = 16
}

foreach ($row in $updates1.Keys) {
    $sheet1.Range("F$row").Value = $updates1[$row]
}

# Row -> new value mapping for sheet "全部类型"
$updates4 = @{
    2  = 14
    3  = 1343
    6  = 60
    8  = 11612
    9  = 4370
    11 = 34
    16 = 1091
    18 = 33
    19 = 5058
    22 = 11327
    23 = 11250
    24 = 16
}

foreach ($row in $updates4.Keys) {
    $sheet4.Range("F$row").Value = $updates4[$row]
}
